$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a new column before D (shifts D:K -> E:L), mirroring adding a new fiscal-year column
$ws.Columns("D:D").Insert()

# Copy number formatting from column E (the old D column, now shifted) into new column D
$ws.Range("E7:E102").Copy()
$ws.Range("D7:D102").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# ---- New column D values (2018 data) ----
$ws.Range("D7").Value = 43465
$ws.Range("D8").Value = 7172200
$ws.Range("D9").Value = 2647500
$ws.Range("D10").Value = 4524700
$ws.Range("D12").Value = "NA"
$ws.Range("D13").Value = 0
$ws.Range("D14").Value = 0
$ws.Range("D15").Value = 0
$ws.Range("D17").Value = 5979900
$ws.Range("D18").Value = 1192300
$ws.Range("D20").Value = 0
$ws.Range("D21").Value = 2860700
$ws.Range("D22").Value = 0
$ws.Range("D23").Value = 1192300
$ws.Range("D24").Value = 276300
$ws.Range("D25").Value = 0
$ws.Range("D26").Value = 915900
$ws.Range("D27").Value = 915900
$ws.Range("D28").Value = 0
$ws.Range("D29").Value = "NA"
$ws.Range("D30").Value = 0
$ws.Range("D31").Value = 0
$ws.Range("D32").Value = 0
$ws.Range("D33").Value = 915900
$ws.Range("D34").Value = 0
$ws.Range("D35").Value = 915900
$ws.Range("D38").Value = 43465
$ws.Range("D41").Value = 148400
$ws.Range("D42").Value = 0
$ws.Range("D43").Value = 410000
$ws.Range("D44").Value = 0
$ws.Range("D45").Value = 29100
$ws.Range("D46").Value = 0
$ws.Range("D47").Value = 26321500
$ws.Range("D48").Value = 14040100
$ws.Range("D49").Value = 109300
$ws.Range("D50").Value = 0
$ws.Range("D51").Value = 0
$ws.Range("D52").Value = 2102000
$ws.Range("D53").Value = 0
$ws.Range("D54").Value = 43959900
$ws.Range("D57").Value = 423000
$ws.Range("D58").Value = 0
$ws.Range("D59").Value = 112600
$ws.Range("D60").Value = 0
$ws.Range("D61").Value = 34883000
$ws.Range("D62").Value = 1155900
$ws.Range("D63").Value = 0
$ws.Range("D64").Value = 0
$ws.Range("D65").Value = 0
$ws.Range("D66").Value = 36941500
$ws.Range("D68").Value = 0
$ws.Range("D69").Value = 0
$ws.Range("D70").Value = 0
$ws.Range("D71").Value = 0
$ws.Range("D72").Value = 5465700
$ws.Range("D73").Value = 0
$ws.Range("D74").Value = 0
$ws.Range("D75").Value = 0
$ws.Range("D76").Value = 7018400
$ws.Range("D77").Value = 0
$ws.Range("D80").Value = 43465
$ws.Range("D81").Value = 915900
$ws.Range("D83").Value = 1668500
$ws.Range("D84").Value = 0
$ws.Range("D85").Value = 0
$ws.Range("D86").Value = 0
$ws.Range("D87").Value = 0
$ws.Range("D88").Value = 0
$ws.Range("D89").Value = 6244900
$ws.Range("D91").Value = -9829800
$ws.Range("D92").Value = 0
$ws.Range("D93").Value = 0
$ws.Range("D94").Value = -10415800
$ws.Range("D96").Value = -180300
$ws.Range("D97").Value = 0
$ws.Range("D98").Value = 0
$ws.Range("D99").Value = 0
$ws.Range("D100").Value = 3339700
$ws.Range("D101").Value = 0
$ws.Range("D102").Value = -831200

# ---- Corrections to shifted cells (restated values) ----
$ws.Range("E10").Value = 4508800
$ws.Range("E17").Value = 5951100
$ws.Range("E18").Value = 804000
$ws.Range("E21").Value = 2207700
$ws.Range("E23").Value = 804000
$ws.Range("E24").Value = 280900
$ws.Range("E26").Value = 523100
$ws.Range("E27").Value = 523100
$ws.Range("E29").Value = 649700
$ws.Range("E33").Value = 1172800
$ws.Range("E35").Value = 1172800
$ws.Range("E43").Value = 424500
$ws.Range("E47").Value = 24772100
$ws.Range("E54").Value = 39402800
$ws.Range("E62").Value = 892400
$ws.Range("E66").Value = 32937100
$ws.Range("E72").Value = 4736300
$ws.Range("E76").Value = 6465700
$ws.Range("E81").Value = 1172800
$ws.Range("E89").Value = 3941300
$ws.Range("E94").Value = -3590300
$ws.Range("F94").Value = -4742400
$ws.Range("E100").Value = -186800
$ws.Range("F100").Value = 931500
$ws.Range("E102").Value = 164200
$ws.Range("F102").Value = 662300

# Autofit the new column to match data width
$ws.Columns("D:D").AutoFit()
